# Update "Sufficient data/forecast_summary_B089FWWN62.xlsx" with corrected
# forecast output:
#   - Forecast Comparison sheet: insert a new "Week_Start_Date" column after
#     "Week", reformat the Week labels (W01 -> W1), refresh the MyForecast
#     numbers, and store is_holiday_week as a boolean.
#   - Summary sheet: refresh the rolling forecast totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date"); everything from the old B
# column onward (ASIN, MyForecast, Amazon forecasts, Product Title,
# is_holiday_week) shifts one column to the right automatically.
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Week_Start_Date"

# Week labels lose their leading zero (W01 -> W1, ... W09 -> W9). W10-W16
# are unaffected but are re-written for clarity/consistency.
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weeks[$i]
}

# Week start dates (written as text, matching the source data).
$ws.Range("B2:B17").NumberFormat = "@"
$weekStarts = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")
for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $weekStarts[$i]
}

# Refreshed MyForecast values (now column D after the insert).
$myForecast = @(31,40,44,45,39,41,41,37,40,41,36,38,38,43,47,46)
for ($i = 0; $i -lt $myForecast.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $myForecast[$i]
}

# is_holiday_week (now column J) becomes a proper boolean column.
$ws.Range("J2:J17").Value = $false

# ---------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9:B11").NumberFormat = "@"
$ws2.Range("B9").Value = "648"
$ws2.Range("B10").Value = "319"
$ws2.Range("B11").Value = "160"
